$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.026.87"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.409.14"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "2.844.17"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "62.082.99"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "2.412.94"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "567.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").Value = "2.527.38"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "0.0₃0939"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "151.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("E41").Value = "  -10.25%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.591"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0918"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  +0.22%  "
